$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "展览" (sheet1): update "想去人数" (column F) counts on several rows.
# ---------------------------------------------------------------------------
$wsExpo = $wb.Worksheets.Item("展览")
$expoUpdates = @{
    2  = 268
    3  = 623
    6  = 2814
    7  = 533
    8  = 53
    10 = 377
    12 = 315
    14 = 5938
    17 = 7
    18 = 106
    21 = 529
    22 = 27
    23 = 23
    24 = 75
    25 = 1304
    27 = 5
    28 = 34
    29 = 2057
    30 = 173
    33 = 3277
}
foreach ($row in $expoUpdates.Keys) {
    $wsExpo.Cells.Item($row, 6).Value = $expoUpdates[$row]
}

# ---------------------------------------------------------------------------
# Sheet "演出" (sheet2): update several F-column counts, then insert a brand
# new row for a newly announced show (2024.04.25, 赵鹏) before the existing
# "2024.04.26 夏川里美" row, shifting the rows below it down by one.
# ---------------------------------------------------------------------------
$wsShow = $wb.Worksheets.Item("演出")
$showUpdates = @{
    7  = 345
    12 = 644
    13 = 106
    23 = 300
    24 = 4050
    25 = 8
    28 = 135
    30 = 69
    33 = 18
}
foreach ($row in $showUpdates.Keys) {
    $wsShow.Cells.Item($row, 6).Value = $showUpdates[$row]
}

# Insert a new row 34 (existing rows 34-36 shift down to 35-37, carrying
# their cell content with them - only the index numbers in column A stay
# tied to row position, so those get re-stamped below).
$wsShow.Cells.Item(34, 1).EntireRow.Insert()

# Match the look of the other index cells in column A (bold, thin border,
# centered/top aligned) since Insert() alone doesn't reliably clone it.
$newIndexCell = $wsShow.Cells.Item(34, 1)
$newIndexCell.Font.Bold = $true
$newIndexCell.HorizontalAlignment = -4108
$newIndexCell.VerticalAlignment = -4160
$newIndexCell.Borders.LineStyle = 1
$newIndexCell.Value = 33

# B34 looks like a bare date ("2024.04.25"); Excel's COM layer would
# normally auto-coerce that into a date serial. Force text so it round
# trips as the literal string, like the sheet's other B-column cells.
$wsShow.Cells.Item(34, 2).NumberFormat = "@"
$wsShow.Cells.Item(34, 2).Value = "2024.04.25"
$wsShow.Cells.Item(34, 2).Style = "Normal"

$wsShow.Cells.Item(34, 3).Value = "上海·赵鹏`"行吟2024·拥抱`"——巡演十周年纪念演唱会"
$wsShow.Cells.Item(34, 4).Value = "丁香路425号(上海科技馆地铁站1号口步行460米) 上海东方艺术中心音乐厅"
$wsShow.Cells.Item(34, 5).Value = "2024.04.25 19:30-04.25 21:00"
$wsShow.Cells.Item(34, 6).Value = 0
$wsShow.Cells.Item(34, 7).Value = 199
$wsShow.Cells.Item(34, 8).Value = "https://show.bilibili.com/platform/detail.html?id=81992"
$wsShow.Cells.Item(34, 9).Value = "//i1.hdslb.com/bfs/openplatform/202402/YgmTp7py1708583402797.jpeg"

# The index column (A) is a fixed "row position - 1" counter, not data that
# travels with the row, so restamp it for the rows that shifted down.
$wsShow.Cells.Item(35, 1).Value = 34
$wsShow.Cells.Item(36, 1).Value = 35
$wsShow.Cells.Item(37, 1).Value = 36

# ---------------------------------------------------------------------------
# Sheet "本地生活" (sheet3): update F-column counts.
# ---------------------------------------------------------------------------
$wsLocal = $wb.Worksheets.Item("本地生活")
$localUpdates = @{
    2  = 1801
    5  = 2592
    8  = 1484
    9  = 415
    12 = 634
}
foreach ($row in $localUpdates.Keys) {
    $wsLocal.Cells.Item($row, 6).Value = $localUpdates[$row]
}

# ---------------------------------------------------------------------------
# Sheet "全部类型" (sheet4): update F-column counts.
# ---------------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("全部类型")
$allUpdates = @{
    2  = 1801
    4  = 2592
    7  = 1484
    8  = 415
    10 = 268
    11 = 623
    12 = 2814
    13 = 53
    14 = 634
    15 = 377
    18 = 315
    20 = 5938
    24 = 106
    27 = 529
    32 = 23
    35 = 300
    36 = 8
    39 = 135
    41 = 34
    42 = 69
    44 = 2057
    47 = 173
    50 = 3277
}
foreach ($row in $allUpdates.Keys) {
    $wsAll.Cells.Item($row, 6).Value = $allUpdates[$row]
}
